# Week 13 logging update for Cowboys Players Data workbook.

$wb = $excel.ActiveWorkbook

$rushing   = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------------------

# E.Elliott (row 4): 1DATT/2DATT updated
$rushing.Range("C4").Value = 103
$rushing.Range("D4").Value = 49

# T.Pollard (row 5): 1DATT/2DATT updated
$rushing.Range("C5").Value = 66
$rushing.Range("D5").Value = 33

# C.Lamb (row 7): 2DATT updated
$rushing.Range("D7").Value = 2

# C.Wilson (row 8): 3DATT updated
$rushing.Range("E8").Value = 0

# New player row: D.Schultz (copy the number formatting used by the rest of
# column A before writing the values)
$rushing.Range("A8").Copy()
$rushing.Range("A9").PasteSpecial(-4122)
$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "D.Schultz"
$rushing.Range("C9").Value = 0
$rushing.Range("D9").Value = 1
$rushing.Range("E9").Value = 0
$rushing.Range("F9").Value = 0

# ---------------------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------------------

# E.Elliott (row 2)
$receiving.Range("C2").Value = 35
$receiving.Range("D2").Value = 28
$receiving.Range("G2").Value = 11
$receiving.Range("H2").Value = 10

# T.Pollard (row 3)
$receiving.Range("C3").Value = 26
$receiving.Range("D3").Value = 22
$receiving.Range("G3").Value = 3
$receiving.Range("H3").Value = 3

# A.Cooper (row 4)
$receiving.Range("C4").Value = 47
$receiving.Range("D4").Value = 35
$receiving.Range("E4").Value = 19
$receiving.Range("F4").Value = 11

# C.Lamb (row 5)
$receiving.Range("C5").Value = 57
$receiving.Range("D5").Value = 38
$receiving.Range("E5").Value = 29
$receiving.Range("F5").Value = 16

# M.Gallup (row 6)
$receiving.Range("C6").Value = 31
$receiving.Range("D6").Value = 26
$receiving.Range("E6").Value = 8
$receiving.Range("G6").Value = 3
$receiving.Range("H6").Value = 3

# N.Brown (row 8)
$receiving.Range("C8").Value = 11
$receiving.Range("D8").Value = 8

# M.Turner (row 9)
$receiving.Range("C9").Value = 8
$receiving.Range("D9").Value = 7

# D.Schultz (row 11)
$receiving.Range("C11").Value = 48
$receiving.Range("D11").Value = 38
$receiving.Range("G11").Value = 6
$receiving.Range("H11").Value = 3

# ---------------------------------------------------------------------------
# Active sheet switches from Receiving back to Rushing
# ---------------------------------------------------------------------------
$rushing.Activate()
